$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.04187148740654
$ws.Range("D2").Value = 1.051201815403074
$ws.Range("E2").Value = 1.039821650265975
$ws.Range("F2").Value = 1.058114308825262
$ws.Range("J2").Value = 1.046950360256459
$ws.Range("K2").Value = 1.053953972737854
$ws.Range("L2").Value = 1.042605748426806
$ws.Range("M2").Value = 1.06084742817026
$ws.Range("N2").Value = 1.018867705991366

$ws.Range("C3").Value = 1.044480530114546
$ws.Range("D3").Value = 1.053794344556942
$ws.Range("E3").Value = 1.042104672355415
$ws.Range("F3").Value = 1.060870762483992
$ws.Range("J3").Value = 1.049197364533617
$ws.Range("K3").Value = 1.056354528623647
$ws.Range("L3").Value = 1.044695193457337
$ws.Range("M3").Value = 1.063412932042585
$ws.Range("N3").Value = 1.019675759831763

$ws.Range("C4").Value = 1.04615936746695
$ws.Range("D4").Value = 1.055462855549147
$ws.Range("E4").Value = 1.043573327982252
$ws.Range("F4").Value = 1.062645297380344
$ws.Range("J4").Value = 1.050642091253775
$ws.Range("K4").Value = 1.057898577006448
$ws.Range("L4").Value = 1.046038290180848
$ws.Range("M4").Value = 1.065063695222774
$ws.Range("N4").Value = 1.020193834311686

$ws.Range("C5").Value = 1.046862965162818
$ws.Range("D5").Value = 1.056162196342182
$ws.Range("E5").Value = 1.044188741922258
$ws.Range("F5").Value = 1.063389202993886
$ws.Range("J5").Value = 1.051247297617153
$ws.Range("K5").Value = 1.058545532444823
$ws.Range("L5").Value = 1.046600844292043
$ws.Range("M5").Value = 1.065755512594657
$ws.Range("N5").Value = 1.020410505283553

$ws.Range("C6").Value = 1.046980975933701
$ws.Range("D6").Value = 1.056279497256248
$ws.Range("E6").Value = 1.044291956351953
$ws.Range("F6").Value = 1.063513986188948
$ws.Range("J6").Value = 1.051348789515655
$ws.Range("K6").Value = 1.058654033945396
$ws.Range("L6").Value = 1.046695178856998
$ws.Range("M6").Value = 1.065871546647892
$ws.Range("N6").Value = 1.020446819728551

$ws.Range("C7").Value = 1.0461687774792
$ws.Range("D7").Value = 1.055472208351387
$ws.Range("E7").Value = 1.043581558996655
$ws.Range("F7").Value = 1.062655245685326
$ws.Range("J7").Value = 1.050650186453369
$ws.Range("K7").Value = 1.057907230077869
$ws.Range("L7").Value = 1.046045815177342
$ws.Range("M7").Value = 1.065072947741051
$ws.Range("N7").Value = 1.020196733884585

$ws.Range("C8").Value = 1.042755207016751
$ws.Range("D8").Value = 1.052079878841492
$ws.Range("E8").Value = 1.040595023663328
$ws.Range("F8").Value = 1.059047782971711
$ws.Range("J8").Value = 1.047711691033278
$ws.Range("K8").Value = 1.054767206369114
$ws.Range("L8").Value = 1.043313762384146
$ws.Range("M8").Value = 1.061716412354042
$ws.Range("N8").Value = 1.019141795176683

$ws.Range("C9").Value = 1.036665337769894
$ws.Range("D9").Value = 1.046030263772808
$ws.Range("E9").Value = 1.035263993507623
$ws.Range("F9").Value = 1.052618565370168
$ws.Range("J9").Value = 1.042460483694759
$ws.Range("K9").Value = 1.049160464999228
$ws.Range("L9").Value = 1.038428974103901
$ws.Range("M9").Value = 1.055727867516467
$ws.Range("N9").Value = 1.017245289056236

$ws.Range("C10").Value = 1.032551179727446
$ws.Range("D10").Value = 1.041944948051616
$ws.Range("E10").Value = 1.031660577821855
$ws.Range("F10").Value = 1.048279616504078
$ws.Range("J10").Value = 1.038906952928113
$ws.Range("K10").Value = 1.045369439635251
$ws.Range("L10").Value = 1.035121752625017
$ws.Range("M10").Value = 1.051681907630734
$ws.Range("N10").Value = 1.01595444007069

$ws.Range("C11").Value = 1.03075592845285
$ws.Range("D11").Value = 1.040162685939146
$ws.Range("E11").Value = 1.030087765637151
$ws.Range("F11").Value = 1.046387344262863
$ws.Range("J11").Value = 1.037354932637966
$ws.Range("K11").Value = 1.043714426013899
$ws.Range("L11").Value = 1.033676926104314
$ws.Range("M11").Value = 1.049916362575612
$ws.Range("N11").Value = 1.015388915221137

$ws.Range("C12").Value = 1.030086940443625
$ws.Range("D12").Value = 1.039498599960275
$ws.Range("E12").Value = 1.029501605844354
$ws.Range("F12").Value = 1.045682363410549
$ws.Range("J12").Value = 1.036776372736509
$ws.Range("K12").Value = 1.043097582477152
$ws.Range("L12").Value = 1.033138268930366
$ws.Range("M12").Value = 1.049258437118932
$ws.Range("N12").Value = 1.015177840786767

$ws.Range("C13").Value = 1.030230539311422
$ws.Range("D13").Value = 1.039641143794279
$ws.Range("E13").Value = 1.029627428339762
$ws.Range("F13").Value = 1.045833680862917
$ws.Range("J13").Value = 1.036900570629133
$ws.Range("K13").Value = 1.043229993638125
$ws.Range("L13").Value = 1.033253903570133
$ws.Range("M13").Value = 1.049399661718107
$ws.Range("N13").Value = 1.015223163232381

$ws.Range("C14").Value = 1.030700674057047
$ws.Range("D14").Value = 1.040107835157791
$ws.Range("E14").Value = 1.030039353609647
$ws.Range("F14").Value = 1.046329113781357
$ws.Range("J14").Value = 1.03730715134173
$ws.Range("K14").Value = 1.043663480760191
$ws.Range("L14").Value = 1.033632441392221
$ws.Range("M14").Value = 1.049862022024539
$ws.Range("N14").Value = 1.015371488560734

$ws.Range("C15").Value = 1.030990052053566
$ws.Range("D15").Value = 1.040395101877317
$ws.Range("E15").Value = 1.030292894265003
$ws.Range("F15").Value = 1.046634084720211
$ws.Range("J15").Value = 1.037557382670715
$ws.Range("K15").Value = 1.043930286279885
$ws.Range("L15").Value = 1.033865406144137
$ws.Range("M15").Value = 1.050146613745548
$ws.Range("N15").Value = 1.015462741658496

$ws.Range("C16").Value = 1.032670027291305
$ws.Range("D16").Value = 1.042062944286951
$ws.Range("E16").Value = 1.031764690868103
$ws.Range("F16").Value = 1.048404909474333
$ws.Range("J16").Value = 1.039009668907738
$ws.Range("K16").Value = 1.045478987321479
$ws.Range("L16").Value = 1.035217366315299
$ws.Range("M16").Value = 1.0517987875594
$ws.Range("N16").Value = 1.015991831371258

$ws.Range("C17").Value = 1.033720082961781
$ws.Range("D17").Value = 1.04310552476503
$ws.Range("E17").Value = 1.032684513438599
$ws.Range("F17").Value = 1.049512035391182
$ws.Range("J17").Value = 1.039917035496743
$ws.Range("K17").Value = 1.046446787790024
$ws.Range("L17").Value = 1.036061949156614
$ws.Range("M17").Value = 1.052831452570964
$ws.Range("N17").Value = 1.016321936643824

$ws.Range("C18").Value = 1.034331237235226
$ws.Range("D18").Value = 1.043712367005878
$ws.Range("E18").Value = 1.033219828073061
$ws.Range("F18").Value = 1.050156508315902
$ws.Range("J18").Value = 1.040445006863655
$ws.Range("K18").Value = 1.047009994367402
$ws.Range("L18").Value = 1.036553351603943
$ws.Range("M18").Value = 1.053432480014797
$ws.Range("N18").Value = 1.016513848429056

$ws.Range("C19").Value = 1.034539402610254
$ws.Range("D19").Value = 1.043919070187846
$ws.Range("E19").Value = 1.033402154631216
$ws.Range("F19").Value = 1.050376039737977
$ws.Range("J19").Value = 1.040624816393088
$ws.Range("K19").Value = 1.04720181578818
$ws.Range("L19").Value = 1.036720700708766
$ws.Range("M19").Value = 1.053637195203283
$ws.Range("N19").Value = 1.016579178761367

$ws.Range("C20").Value = 1.033607559590002
$ws.Range("D20").Value = 1.042993798418151
$ws.Range("E20").Value = 1.032585950028396
$ws.Range("F20").Value = 1.049393385779098
$ws.Range("J20").Value = 1.039819816565332
$ws.Range("K20").Value = 1.046343086438366
$ws.Range("L20").Value = 1.035971460935525
$ws.Range("M20").Value = 1.052720793296082
$ws.Range("N20").Value = 1.016286585128486

$ws.Range("C21").Value = 1.030562291138301
$ws.Range("D21").Value = 1.039970464078726
$ws.Range("E21").Value = 1.029918106193661
$ws.Range("F21").Value = 1.046183279955624
$ws.Range("J21").Value = 1.037187481177587
$ws.Range("K21").Value = 1.043535888114372
$ws.Range("L21").Value = 1.033521026740186
$ws.Range("M21").Value = 1.049725927686526
$ws.Range("N21").Value = 1.015327838624566

$ws.Range("C22").Value = 1.028635119882757
$ws.Range("D22").Value = 1.038057531862724
$ws.Range("E22").Value = 1.028229425464244
$ws.Range("F22").Value = 1.04415272847382
$ws.Range("J22").Value = 1.035520412689214
$ws.Range("K22").Value = 1.041758716308536
$ws.Range("L22").Value = 1.031968827991334
$ws.Range("M22").Value = 1.047830611641984
$ws.Range("N22").Value = 1.0147191610485

$ws.Range("C23").Value = 1.029657960769318
$ws.Range("D23").Value = 1.039072781158901
$ws.Range("E23").Value = 1.02912572142676
$ws.Range("F23").Value = 1.045230349335697
$ws.Range("J23").Value = 1.036405319775499
$ws.Range("K23").Value = 1.04270200770965
$ws.Range("L23").Value = 1.032792791232976
$ws.Range("M23").Value = 1.048836549195752
$ws.Range("N23").Value = 1.01504239798001

$ws.Range("C24").Value = 1.033658408149573
$ws.Range("D24").Value = 1.043044286687923
$ws.Range("E24").Value = 1.032630490302435
$ws.Range("F24").Value = 1.049447002428488
$ws.Range("J24").Value = 1.039863749567876
$ws.Range("K24").Value = 1.046389948613271
$ws.Range("L24").Value = 1.036012352460081
$ws.Range("M24").Value = 1.052770799498561
$ws.Range("N24").Value = 1.016302560911395

$ws.Range("C25").Value = 1.03824897671035
$ws.Range("D25").Value = 1.047603150431547
$ws.Range("E25").Value = 1.036650642989912
$ws.Range("F25").Value = 1.054289671377581
$ws.Range("J25").Value = 1.043827077826951
$ws.Range("K25").Value = 1.050619041999529
$ws.Range("L25").Value = 1.039700501341906
$ws.Range("M25").Value = 1.057285204556185
$ws.Range("N25").Value = 1.017740160577004
